$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("locations")

$ws.Range("A2").Value = "4fizzqrs626"
$ws.Range("B2").Value = "foidgiuw"
$ws.Range("C2").Value = "HKD"
$ws.Range("D2").Value = "Madam Thu Bakery, 21C, Võ Văn Tần, Ninh Kiều, Ninh Kiều District, Cần Thơ, 94111, Vietnam"
$ws.Range("E2").Value = "https://www.google.com/maps/search/?api=1&query=10.032100,105.786400"
$ws.Range("F2").Value = "2025-08-22T06:46:56.105Z"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = "Nguyễn Văn B"
$ws.Range("M2").Value = "150fea0e0fdf86f5"
$ws.Range("N2").Value = "178f60f6166cdc837bfab2bbc150a80dd6fc48d467309288fc27afb78a08279c"
